# Weekly update: a new "Choclero" record for Vega Monumental Concepción is
# inserted at the top of the data block (row 17), pushing the existing
# records (previously rows 17-61) down by one row (now rows 18-62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 17; existing rows 17..61 shift to 18..62.
$ws.Rows("17").Insert()

# Populate the new row 17 with the newly added record.
$ws.Cells.Item(17, 1).Value2  = 11
$ws.Cells.Item(17, 2).Value2  = "Vega Monumental Concepción"
$ws.Cells.Item(17, 3).Value2  = "Bíobío"
$ws.Cells.Item(17, 4).Value2  = 44560
$ws.Cells.Item(17, 5).Value2  = 8
$ws.Cells.Item(17, 6).Value2  = 100112024
$ws.Cells.Item(17, 7).Value2  = "Choclo"
$ws.Cells.Item(17, 8).Value2  = "Choclero"
$ws.Cells.Item(17, 9).Value2  = "Primera"
$ws.Cells.Item(17, 10).Value2 = 2000
$ws.Cells.Item(17, 11).Value2 = 300
$ws.Cells.Item(17, 12).Value2 = 350
$ws.Cells.Item(17, 13).Value2 = 325
$ws.Cells.Item(17, 14).Value2 = "$/unidad"
$ws.Cells.Item(17, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(17, 16).Value2 = 325
$ws.Cells.Item(17, 17).Value2 = 1
$ws.Cells.Item(17, 18).Value2 = "Hortaliza"
